# Hotfix: Wed Nov 20 17:48:22 RTZ 2024
$wb = $excel.ActiveWorkbook

# --- SQL sheet: row 13 ---
$wsSQL = $wb.Worksheets.Item("SQL")
$wsSQL.Range("B13").Value = 'Создание таблицы в базе данных.'
$wsSQL.Range("E13").Value = '2024-11-20 12:29:35'

# --- Python sheet: row 36 ---
$wsPython = $wb.Worksheets.Item("Python")
$wsPython.Range("C36").Value = 'Вывод текущей даты и времени.'
$wsPython.Range("E36").Value = '2024-11-20 12:26:54'

# --- HTML sheet: add new row 1 ---
$wsHTML = $wb.Worksheets.Item("HTML")
$wsHTML.Range("A1").Value = 3
$wsHTML.Range("B1").Value = 'Тег h1'
$wsHTML.Range("C1").Value = 'Тег h1 - заголовок первого уровня. Пример написания: <h1></h1>'
$wsHTML.Range("D1").Value = '2024-11-20 09:36:48'
$wsHTML.Range("E1").Value = '2024-11-20 13:47:03'

# --- Bash sheet: rows 80-81 ---
$wsBash = $wb.Worksheets.Item("Bash")
$wsBash.Range("C80").Value = 'Скрипт powershell, который убивает все python запущенные процессы.'
$wsBash.Range("E80").Value = '2024-11-20 16:21:09'
$wsBash.Range("C81").Value = 'Проверка обновлений библиотек!!'
$wsBash.Range("E81").Value = '2024-11-20 12:36:21'

# --- Test sheet: reshuffled option values, drop trailing 2 rows ---
$wsTest = $wb.Worksheets.Item("Test")

$testVals = @(
    'Samsung',
    'AMD',
    'ASUS',
    'Acer',
    'ASUS',
    'ASUS',
    'ASUS',
    'Samsung',
    'Toyota',
    'ASUS',
    'HP',
    'Ford',
    'Intel',
    'AMD',
    'ASUS',
    'Acer',
    'Toyota',
    'HP',
    'Toyota',
    'HP',
    'Acer',
    'Ford',
    'Acer',
    'AMD',
    'Ford',
    'HP',
    'Intel',
    'HP',
    'Acer',
    'HP',
    'ASUS',
    'Acer',
    'Samsung',
    'Toyota',
    'Toyota',
    'AMD',
    'Acer',
    'Ford',
    'HP',
    'Intel',
    'Ford',
    'Acer',
    'ASUS',
    'Samsung',
    'HP',
    'Intel',
    'Samsung',
    'Intel',
    'Acer',
    'Ford',
    'Acer',
    'Acer',
    'AMD',
    'Samsung',
    'Samsung',
    'Ford',
    'Acer',
    'Acer',
    'Intel',
    'Ford',
    'Toyota',
    'Ford',
    'Samsung',
    'HP',
    'Intel',
    'Intel',
    'Toyota',
    'AMD',
    'AMD',
    'Intel',
    'Acer',
    'Toyota',
    'Intel',
    'Intel',
    'Ford',
    'ASUS',
    'Intel',
    'Ford',
    'ASUS',
    'ASUS',
    'AMD',
    'Acer',
    'HP',
    'HP',
    'AMD',
    'ASUS',
    'AMD',
    'AMD',
    'Acer',
    'Samsung',
    'Toyota',
    'Ford',
    'Samsung',
    'HP',
    'ASUS',
    'ASUS',
    'Acer',
    'Toyota',
    'Toyota',
    'Intel',
    'Тестовая запись №1',
    'Тестовая запись №2',
    'Тестовая запись №3',
    'Тестовая запись №4',
    'Тестовая запись №5',
    'Тестовая запись №6',
    'Тестовая запись №7',
    'Тестовая запись №8',
    'Тестовая запись №9',
    'Тестовая запись №10',
    'Тестовая запись №11',
    'Тестовая запись №12',
    'Тестовая запись №13',
    'Тестовая запись №14',
    'Тестовая запись №15',
    'Тестовая запись №16',
    'Тестовая запись №17',
    'Тестовая запись №18',
    'Тестовая запись №19',
    'Тестовая запись №20',
    'option3'
)

for ($i = 0; $i -lt $testVals.Length; $i++) {
    $wsTest.Cells.Item($i + 1, 2).Value = $testVals[$i]
}

$wsTest.Rows("122:123").Delete()
